# Add a new "2023_TM160_IPA_40" model run (AOC = 16.46 cents) to the
# ModelRuns_RTP2025 log, inserted right after the current last 2023 run
# (row 51, "2023_TM160_IPA_39") and before the first 2025 run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 52; this shifts the existing rows 52-60 down to 53-61
# and (in this engine) automatically inherits the formatting/styles of the
# row above it (row 51), which already matches what the new row needs.
$ws.Rows(52).Insert()

# Fill in the new row's values (columns A-S).
$ws.Cells.Item(52, 1).Value  = 2023
$ws.Cells.Item(52, 2).Value  = "2023_TM160_IPA_40"
$ws.Cells.Item(52, 3).Value  = "RTP2025_IP"
$ws.Cells.Item(52, 4).Value  = "Past year"
$ws.Cells.Item(52, 5).Value  = "new AOC (16.46 cents)"
$ws.Cells.Item(52, 6).Value  = "petrale"
$ws.Cells.Item(52, 7).Value  = "n/a"
$ws.Cells.Item(52, 8).Value  = "current"
$ws.Cells.Item(52, 9).Value  = "BlueprintNetworks_v12\net_2023_Blueprint"
$ws.Cells.Item(52, 10).Value = "model3-c"
$ws.Cells.Item(52, 11).Value = "https://app.asana.com/0/1204085012544660/1205973396667333/f"
$ws.Cells.Item(52, 12).Value = 16.45
$ws.Cells.Item(52, 13).Value = "na"
$ws.Cells.Item(52, 14).Value = "na"
$ws.Cells.Item(52, 15).Value = 0.99
$ws.Cells.Item(52, 16).Value = 0.89
$ws.Cells.Item(52, 17).Value = 100
$ws.Cells.Item(52, 18).Value = 0
$ws.Cells.Item(52, 19).Value = 75

# Update the window view: scroll the frozen (bottom-right) pane and move the
# selection to reflect where the author left off editing.
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 34
[void]$ws.Range("T52").Select()
